$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]1.120168333333333
$ws.Cells.Item(2, 8).Value = [double]3.360505
$ws.Cells.Item(2, 9).Value = [double]0.001768092629909379
$ws.Cells.Item(2, 10).Value = [double]0.001768092629909379
$ws.Cells.Item(2, 13).Value = [double]2.582049666666667
$ws.Cells.Item(2, 14).Value = [double]7.746149
$ws.Cells.Item(2, 15).Value = [double]0.02216861605835153
$ws.Cells.Item(2, 16).Value = [double]0.02216861605835152
$ws.Cells.Item(2, 17).Value = [double]2.892330271693889
$ws.Cells.Item(2, 18).Value = [double]26.030972445245
$ws.Cells.Item(2, 19).Value = [double]0.00003919616666806203
$ws.Cells.Item(2, 20).Value = [double]0.00003919616666806202
$ws.Cells.Item(3, 7).Value = [double]1.120168333333333
$ws.Cells.Item(3, 8).Value = [double]3.360505
$ws.Cells.Item(3, 9).Value = [double]0.001768092629909379
$ws.Cells.Item(3, 10).Value = [double]0.001768092629909379
$ws.Cells.Item(3, 13).Value = [double]93.97803500000002
$ws.Cells.Item(3, 15).Value = [double]0.8068640207540504
$ws.Cells.Item(3, 16).Value = [double]0.8068640207540503
$ws.Cells.Item(3, 17).Value = [double]105.2712188358917
$ws.Cells.Item(3, 18).Value = [double]947.4409695230252
$ws.Cells.Item(3, 19).Value = [double]0.001426610328434284
$ws.Cells.Item(3, 20).Value = [double]0.001426610328434284
$ws.Cells.Item(4, 7).Value = [double]1.120168333333333
$ws.Cells.Item(4, 8).Value = [double]3.360505
$ws.Cells.Item(4, 9).Value = [double]0.001768092629909379
$ws.Cells.Item(4, 10).Value = [double]0.001768092629909379
$ws.Cells.Item(4, 14).Value = [double]59.73934799999999
$ws.Cells.Item(4, 15).Value = [double]0.1709673631875981
$ws.Cells.Item(4, 16).Value = [double]0.1709673631875981
$ws.Cells.Item(4, 17).Value = [double]22.30604196119333
$ws.Cells.Item(4, 18).Value = [double]200.75437765074
$ws.Cells.Item(4, 19).Value = [double]0.0003022861348070322
$ws.Cells.Item(4, 20).Value = [double]0.0003022861348070322
$ws.Cells.Item(5, 9).Value = [double]0.9534130698726969
$ws.Cells.Item(5, 10).Value = [double]0.9534130698726969
$ws.Cells.Item(5, 13).Value = [double]2.582049666666667
$ws.Cells.Item(5, 14).Value = [double]7.746149
$ws.Cells.Item(5, 15).Value = [double]0.02216861605835153
$ws.Cells.Item(5, 16).Value = [double]0.02216861605835152
$ws.Cells.Item(5, 17).Value = [double]1559.638582715397
$ws.Cells.Item(5, 18).Value = [double]14036.74724443857
$ws.Cells.Item(5, 19).Value = [double]0.02113584829102209
$ws.Cells.Item(5, 20).Value = [double]0.02113584829102209
$ws.Cells.Item(6, 9).Value = [double]0.9534130698726969
$ws.Cells.Item(6, 10).Value = [double]0.9534130698726969
$ws.Cells.Item(6, 13).Value = [double]93.97803500000002
$ws.Cells.Item(6, 15).Value = [double]0.8068640207540504
$ws.Cells.Item(6, 16).Value = [double]0.8068640207540503
$ws.Cells.Item(6, 17).Value = [double]56765.66613182034
$ws.Cells.Item(6, 19).Value = [double]0.7692747029969467
$ws.Cells.Item(6, 20).Value = [double]0.7692747029969466
$ws.Cells.Item(7, 9).Value = [double]0.9534130698726969
$ws.Cells.Item(7, 10).Value = [double]0.9534130698726969
$ws.Cells.Item(7, 14).Value = [double]59.73934799999999
$ws.Cells.Item(7, 15).Value = [double]0.1709673631875981
$ws.Cells.Item(7, 16).Value = [double]0.1709673631875981
$ws.Cells.Item(7, 18).Value = [double]108253.2918516745
$ws.Cells.Item(7, 19).Value = [double]0.1630025185847282
$ws.Cells.Item(7, 20).Value = [double]0.1630025185847282
$ws.Cells.Item(8, 8).Value = [double]85.18441
$ws.Cells.Item(8, 9).Value = [double]0.04481883749739363
$ws.Cells.Item(8, 10).Value = [double]0.04481883749739363
$ws.Cells.Item(8, 13).Value = [double]2.582049666666667
$ws.Cells.Item(8, 14).Value = [double]7.746149
$ws.Cells.Item(8, 15).Value = [double]0.02216861605835153
$ws.Cells.Item(8, 16).Value = [double]0.02216861605835152
$ws.Cells.Item(8, 17).Value = [double]73.31679248189889
$ws.Cells.Item(8, 18).Value = [double]659.85113233709
$ws.Cells.Item(8, 19).Value = [double]0.0009935716006613678
$ws.Cells.Item(8, 20).Value = [double]0.0009935716006613678
$ws.Cells.Item(9, 8).Value = [double]85.18441
$ws.Cells.Item(9, 9).Value = [double]0.04481883749739363
$ws.Cells.Item(9, 10).Value = [double]0.04481883749739363
$ws.Cells.Item(9, 13).Value = [double]93.97803500000002
$ws.Cells.Item(9, 15).Value = [double]0.8068640207540504
$ws.Cells.Item(9, 16).Value = [double]0.8068640207540503
$ws.Cells.Item(9, 17).Value = [double]2668.487821478117
$ws.Cells.Item(9, 19).Value = [double]0.03616270742866942
$ws.Cells.Item(9, 20).Value = [double]0.03616270742866941
$ws.Cells.Item(10, 8).Value = [double]85.18441
$ws.Cells.Item(10, 9).Value = [double]0.04481883749739363
$ws.Cells.Item(10, 10).Value = [double]0.04481883749739363
$ws.Cells.Item(10, 14).Value = [double]59.73934799999999
$ws.Cells.Item(10, 15).Value = [double]0.1709673631875981
$ws.Cells.Item(10, 16).Value = [double]0.1709673631875981
$ws.Cells.Item(10, 17).Value = [double]565.4290125738532
$ws.Cells.Item(10, 18).Value = [double]5088.861113164679
$ws.Cells.Item(10, 19).Value = [double]0.007662558468062837
$ws.Cells.Item(10, 20).Value = [double]0.007662558468062837
